# Apply updated cryptocurrency price/volume data scraped on Mon May  8 06:38:56 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D, E
$ws.Range("D2").Value = "28.319.63"
$ws.Range("E2").Value = "  -2.31%  "

# Row 3: update D, E
$ws.Range("D3").Value = "1.870.12"
$ws.Range("E3").Value = "  -1.99%  "

# Row 4: update D, E
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5: update D, E
$ws.Range("D5").Value = "'318.82"
$ws.Range("E5").Value = "  -1.79%  "

# Row 6: update D, E
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.03%  "

# Row 7: update D, E
$ws.Range("D7").Value = "'0.4413"
$ws.Range("E7").Value = "  -3.97%  "

# Row 8: update D
$ws.Range("D8").Value = "'0.3701"

# Row 9: update D, E
$ws.Range("D9").Value = "'0.07507"
$ws.Range("E9").Value = "  -2.68%  "

# Row 10: update D, E
$ws.Range("D10").Value = "'0.9384"
$ws.Range("E10").Value = "  -4.31%  "

# Row 11: update D, E
$ws.Range("D11").Value = "'21.46"
$ws.Range("E11").Value = "  -2.68%  "

# Row 12: update D, E
$ws.Range("D12").Value = "1.860.15"
$ws.Range("E12").Value = "  -2.49%  "

# Row 13: update D, E
$ws.Range("D13").Value = "'6.712"
$ws.Range("E13").Value = "  -3.19%  "

# Row 14: update D, E
$ws.Range("D14").Value = "'5.471"
$ws.Range("E14").Value = "  -3.48%  "

# Row 15: update D, E
$ws.Range("D15").Value = "'0.06883"
$ws.Range("E15").Value = "  -2.53%  "

# Row 16: update E
$ws.Range("E16").Value = "  -0.17%  "

# Row 17: update D, E
$ws.Range("D17").Value = "'82.13"
$ws.Range("E17").Value = "  -1.87%  "

# Row 18: update D, E
$ws.Range("D18").Value = "'0.000009051"
$ws.Range("E18").Value = "  -4.36%  "

# Row 20: update D, E
$ws.Range("D20").Value = "'15.93"
$ws.Range("E20").Value = "  -4.36%  "

# Row 21: update D, E
$ws.Range("D21").Value = "28.302.74"
$ws.Range("E21").Value = "  -2.25%  "

# Row 22: update D, E
$ws.Range("D22").Value = "'5.123"
$ws.Range("E22").Value = "  -3.61%  "

# Row 23: update D, E
$ws.Range("D23").Value = "'10.86"
$ws.Range("E23").Value = "  -0.19%  "

# Row 24: update D, E
$ws.Range("D24").Value = "2.131.26"
$ws.Range("E24").Value = "  -0.15%  "

# Row 25: update D, E
$ws.Range("D25").Value = "'2.029"
$ws.Range("E25").Value = "  -3.11%  "

# Row 26: update D, E
$ws.Range("D26").Value = "'154.94"
$ws.Range("E26").Value = "  -2.11%  "

# Row 27: update D, E
$ws.Range("D27").Value = "'18.38"
$ws.Range("E27").Value = "  -3.71%  "

# Row 28: update D, E
$ws.Range("D28").Value = "'5.336"
$ws.Range("E28").Value = "  -5.53%  "

# Row 29: update D, E
$ws.Range("D29").Value = "'113.73"
$ws.Range("E29").Value = "  -3.18%  "

# Row 30: update D, E
$ws.Range("D30").Value = "'1.725"
$ws.Range("E30").Value = "  -6.91%  "

# Row 31: update D, E
$ws.Range("D31").Value = "'0.09041"
$ws.Range("E31").Value = "  -2.74%  "

# Row 32: update D, E
$ws.Range("D32").Value = "'0.7996"
$ws.Range("E32").Value = "  -7.19%  "

# Row 33: update D, E
$ws.Range("D33").Value = "'4.864"
$ws.Range("E33").Value = "  -4.17%  "

# Row 34: update D, E
$ws.Range("D34").Value = "'1.178"
$ws.Range("E34").Value = "  -5.28%  "

# Row 35: update D, E
$ws.Range("D35").Value = "'2.924"
$ws.Range("E35").Value = "  -1.93%  "

# Row 36: update B, C, D, E
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.131"
$ws.Range("E36").Value = "  -1.90%  "

# Row 37: update B, C, D, E
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  -0.06%  "

# Row 38: update D, E
$ws.Range("D38").Value = "'0.05450"
$ws.Range("E38").Value = "  -4.81%  "

# Row 39: update B, C, D, E
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'3.060"
$ws.Range("E39").Value = "  +7.61%  "

# Row 40: update B, C, D, E
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01969"
$ws.Range("E40").Value = "  -3.30%  "

# Row 41: update D, E
$ws.Range("D41").Value = "'7.112"
$ws.Range("E41").Value = "  -3.82%  "

# Row 42: update D, E
$ws.Range("D42").Value = "'0.5270"
$ws.Range("E42").Value = "  -4.28%  "

# Row 43: update D, E
$ws.Range("D43").Value = "'0.1686"
$ws.Range("E43").Value = "  -3.97%  "

# Row 44: update D, E
$ws.Range("D44").Value = "'8.752"
$ws.Range("E44").Value = "  -6.17%  "

# Row 45: update D, E
$ws.Range("D45").Value = "'0.06762"
$ws.Range("E45").Value = "  -0.91%  "

# Row 46: update D, E
$ws.Range("D46").Value = "'0.4883"
$ws.Range("E46").Value = "  -5.80%  "

# Row 47: update D, E
$ws.Range("D47").Value = "'10.62"
$ws.Range("E47").Value = "  -5.52%  "

# Row 48: update D, E
$ws.Range("D48").Value = "'107.68"
$ws.Range("E48").Value = "  -3.04%  "

# Row 49: update D, E
$ws.Range("D49").Value = "'1.959"
$ws.Range("E49").Value = "  -4.46%  "

# Row 50: update B, C, D, E
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.678"
$ws.Range("E50").Value = "  -5.52%  "

# Row 51: update B, C, D, E
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  -0.13%  "
